$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A78").Value = 45983
$ws.Range("A78").Style = $ws.Range("A77").Style
$ws.Range("A78").NumberFormat = $ws.Range("A77").NumberFormat

$ws.Range("B78").Value = "21,4507"
$ws.Range("C78").Value = "15,7543"
$ws.Range("D78").Value = "15,2104"
$ws.Range("E78").Value = "15,2104"
